$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid Excel auto-converting numeric-looking strings
# like "606.64" into actual numbers) for the Price/Volume columns being touched,
# then clear the temporary formatting so cell style stays the default (as before).
$touchedRange = $ws.Range("D2:E51")
$touchedRange.NumberFormat = "@"

$ws.Range("D2").Value = "71.225.78"
$ws.Range("E2").Value = "  +4.27%  "
$ws.Range("D3").Value = "2.624.55"
$ws.Range("E3").Value = "  +4.55%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "606.64"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "180.99"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "2.623.52"
$ws.Range("E9").Value = "  +4.58%  "
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +15.20%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("D13").Value = "5.04"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "3.085.33"
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  +8.87%  "
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "71.189.21"
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").Value = "2.621.90"
$ws.Range("E18").Value = "  +5.03%  "
$ws.Range("D19").Value = "383.51"
$ws.Range("E19").Value = "  +9.20%  "
$ws.Range("D20").Value = "7.91"
$ws.Range("E20").Value = "  +6.32%  "
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "72.20"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("E24").Value = "  +6.19%  "
$ws.Range("E26").Value = "  +9.51%  "
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  +5.48%  "
$ws.Range("D28").Value = "2.758.74"
$ws.Range("E28").Value = "  +4.42%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "0.0₃0963"
$ws.Range("E30").Value = "  +7.57%  "
$ws.Range("D31").Value = "544.75"
$ws.Range("E31").Value = "  +6.60%  "
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "164.43"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").Value = "19.18"
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  +7.99%  "
$ws.Range("D40").Value = "19.02"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").Value = "2.67"
$ws.Range("E42").Value = "  +9.96%  "
$ws.Range("E43").Value = "  +5.32%  "
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "40.15"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").Value = "154.83"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").Value = "1.70"
$ws.Range("E49").Value = "  +6.18%  "
$ws.Range("D50").Value = "0.533"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  +2.48%  "

$touchedRange.ClearFormats()
